$wb = $excel.ActiveWorkbook

# Update both "展览" and "全部类型" sheets with the revised "想去人数" (column F) counts.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 258
    $ws.Range("F7").Value = 7000
    $ws.Range("F16").Value = 24
    $ws.Range("F18").Value = 633
    $ws.Range("F19").Value = 13
}
